$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.782.55'
$ws.Range("E2").Value = '  +1.36%  '
$ws.Range("D3").Value = '3.169.85'
$ws.Range("E3").Value = '  +1.27%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '617.45'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +1.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.46'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  -2.15%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.167.40'
$ws.Range("E8").Value = '  +1.19%  '
$ws.Range("E9").Value = '  -0.43%  '
$ws.Range("E10").Value = '  -0.03%  '
$ws.Range("E11").Value = '  -1.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.476'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000262'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  +1.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.96'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  -2.51%  '
$ws.Range("D15").Value = '3.685.58'
$ws.Range("E15").Value = '  +2.25%  '
$ws.Range("E16").Value = '  +3.14%  '
$ws.Range("D17").Value = '64.761.73'
$ws.Range("E17").Value = '  +1.16%  '
$ws.Range("D18").Value = '3.165.55'
$ws.Range("E18").Value = '  +1.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.95'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  -0.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '481.64'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  -0.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.78'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +0.79%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.721'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  +1.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.00'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +2.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.85'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.73'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  +0.55%  '
$ws.Range("E26").Value = '  -0.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.84'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  -3.58%  '
$ws.Range("E28").Value = '  +0.71%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.119'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  -5.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.94'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -1.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.11'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  -6.64%  '
$ws.Range("E32").Value = '  +0.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.72'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  +0.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.76'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("E35").Value = '  +2.73%  '
$ws.Range("D36").Value = '0.0₃0791'
$ws.Range("E36").Value = '  +5.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.05'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  -1.08%  '
$ws.Range("E38").Value = '  -1.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '53.19'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  -2.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '466.63'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  +3.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0402'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  +0.37%  '
$ws.Range("E42").Value = '  -3.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.42'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -1.18%  '
$ws.Range("D44").Value = '2.854.54'
$ws.Range("E44").Value = '  -1.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.35'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  +0.63%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.270'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -0.99%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.46'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  +5.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '26.84'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  +0.21%  '
$ws.Range("E49").Value = '  +0.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.115'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  -0.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '35.36'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  +6.53%  '
